# Apply gradebook corrections to "Лист1" (sheet 1).
# Column E/F/G/D hold raw scores that feed the J-column MAX(...) formula,
# which recalculates automatically. Column N holds free-text notes for a
# few students (merged N:S per row), referencing existing / new shared
# strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# New note strings are appended to the shared-strings table in the order
# they are first written, so seed them in the same order as the target
# workbook: "пропуск по уважительной причине (военкомат)" (row 26) first,
# then "деление на 0 не обработано" (rows 7 & 21).
$ws.Range("N26").Value2 = "пропуск по уважительной причине (военкомат)"
$ws.Range("N7").Value2 = "деление на 0 не обработано"

# --- Row 3: Ардаширова Амина Рифовна ---
$ws.Range("E3").Value2 = 4

# --- Row 4: Арзамаскина Юлиана Анатольевна ---
$ws.Range("F4").Value2 = 4
$ws.Range("G4").Value2 = 0
$ws.Range("N4").Value2 = "переписаны верно все номера"

# --- Row 7: Галямова Яна Дмитриевна ---
$ws.Range("F7").Value2 = 4

# --- Row 9: Жамсаранова Аяна Жаргаловна ---
$ws.Range("E9").Value2 = 4

# --- Row 10: Захаренкова Екатерина Денисовна ---
$ws.Range("F10").Value2 = 4
$ws.Range("N10").Value2 = "переписаны верно все номера"

# --- Row 11: Иванов Дмитрий Сергеевич ---
$ws.Range("E11").Value2 = 4

# --- Row 14: Круглов Кирилл Максимович ---
$ws.Range("E14").Value2 = 5

# --- Row 18: Мачкалян Тигран Норайрович ---
$ws.Range("E18").Value2 = 4

# --- Row 19: Молокова Татьяна Михайловна ---
$ws.Range("E19").Value2 = 4

# --- Row 20: Нефодина Анна Александровна ---
$ws.Range("F20").Value2 = 4
$ws.Range("N20").Value2 = "переписаны верно все номера"

# --- Row 21: Ротанкова Вера Владимировна ---
$ws.Range("E21").Value2 = 4
$ws.Range("F21").Value2 = 3
$ws.Range("N21").Value2 = "деление на 0 не обработано"

# --- Row 22: Сахно Полина Валерьевна ---
$ws.Range("D22").Value2 = 4

# --- Row 26: Шаблыгин Михаил Максимович ---
$ws.Range("F26").Value2 = 5

$wb.Save()
